$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set number format to Text ("@") for D-column price cells whose new values
# would otherwise be auto-parsed by Excel as numbers, so they remain text
# (matching the original inlineStr/shared-string text cells).
$ws.Range('D2').Value = '37.277.64'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').Value = '2.002.03'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.14'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.630'
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.57'
$ws.Range('E7').Value = '  +4.20%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.384'
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0803'
$ws.Range('E10').Value = '  +0.48%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.89'
$ws.Range('E12').Value = '  +9.25%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.73'
$ws.Range('E13').Value = '  +4.58%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.849'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('D15').Value = '2.299.60'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.44'
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').Value = '2.001.35'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('D18').Value = '37.181.33'
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.48'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('D20').Value = '0.0₃0867'
$ws.Range('E20').Value = '  +1.82%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.22'
$ws.Range('E21').Value = '  +3.16%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '231.44'
$ws.Range('E22').Value = '  +1.51%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.52'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.42'
$ws.Range('E26').Value = '  +2.46%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.141'
$ws.Range('E27').Value = '  +2.60%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '164.02'
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.74'
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.35'
$ws.Range('E30').Value = '  +17.55%  '
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.86'
$ws.Range('E32').Value = '  +3.37%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0643'
$ws.Range('E33').Value = '  +4.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.56'
$ws.Range('E34').Value = '  +5.59%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.38'
$ws.Range('E35').Value = '  +5.78%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  +2.57%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.34'
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.51'
$ws.Range('E39').Value = '  -3.83%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0977'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.95'
$ws.Range('E41').Value = '  +1.75%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.18'
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0214'
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.70'
$ws.Range('E44').Value = '  +5.73%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '90.97'
$ws.Range('E45').Value = '  +3.84%  '
$ws.Range('D46').Value = '1.377.88'
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('E47').Value = '  +2.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.28'
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.03'
$ws.Range('E49').Value = '  +14.83%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.85'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '46.62'
$ws.Range('E51').Value = '  +7.52%  '
